# POM Feb 4 th
$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestSteps = $wb.Worksheets.Item("TestSteps")

# Mark the "Role updation" test case (TC005, row 4) as Executed = Y
$wsTestCases.Range("C4").Value = "Y"

# Update the remembered selection on the TestSteps sheet
$wsTestSteps.Range("E2").Select()

# Update the remembered selection on the TestCases sheet and make it the
# active sheet/tab (selecting a range activates its parent sheet, so this
# must run last so TestCases ends up as the selected tab)
$wsTestCases.Range("C4").Select()
